# Auto-update predictions and index for 2025-10-30
#
# Populates Sheet1 with the fixtures/predictions table: a bold, boxed,
# centered/top-aligned header row followed by five data rows of
# fixture / pick / confidence / odds / result values. A few cells in the
# confidence/result columns have no observation yet for this slate, so
# they are written as blank (but present) text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell that never receives any value/formatting -- used purely as a
# "blank format" source so blank-but-present cells don't pick up stray
# style bits (see Set-BlankCell below).
$blankFormatSource = $ws.Cells.Item(200, 26)

function Set-BlankCell($row, $col) {
    # A leading single-quote is Excel's "force text" prefix; assigning it
    # alone stores an empty text value in the cell instead of clearing it
    # outright (plain "" clears/removes the cell entirely). Re-pasting
    # formats from a never-touched cell afterwards strips the quote-prefix
    # formatting flag that the assignment leaves behind, so the cell ends
    # up as a plain, unstyled, empty text cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'"
    $blankFormatSource.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

# ---- Header row (A1:G1) ----------------------------------------------
$headers = @("Fixture", "Pick", "AI_Confidence", "OLBG_Confidence", "Oddspedia_Confidence", "Odds", "Result")

$firstHeader = $ws.Cells.Item(1, 1)
$firstHeader.Value = $headers[0]
$firstHeader.HorizontalAlignment = -4108   # xlCenter
$firstHeader.VerticalAlignment = -4160     # xlTop
$firstHeader.Font.Bold = $true
$firstHeader.Borders.LineStyle = 1         # xlContinuous (thin box border)

# Copy the fully-built header format onto the rest of the header cells so
# they all share the same single style entry.
$firstHeader.Copy()
for ($col = 2; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.PasteSpecial(-4122)              # xlPasteFormats
}

# ---- Data rows (A2:G6) -------------------------------------------------

# Row 2 -- Ajman Club - Al-Ain FC ✓: 0:3
$ws.Cells.Item(2, 1).Value = "Ajman Club - Al-Ain FC ✓: 0:3"
$ws.Cells.Item(2, 2).Value = "Al-Ain FC"
$ws.Cells.Item(2, 3).Value = 38
$ws.Cells.Item(2, 4).Value = 40
Set-BlankCell 2 5
$ws.Cells.Item(2, 6).Value = 5
$ws.Cells.Item(2, 7).Value = "✓"

# Row 3 -- Grasshopper Club Zurich - BSC Young Boys
$ws.Cells.Item(3, 1).Value = "Grasshopper Club Zurich - BSC Young Boys : -:-'"
$ws.Cells.Item(3, 2).Value = "BSC Young Boys"
$ws.Cells.Item(3, 3).Value = 30
$ws.Cells.Item(3, 4).Value = 54
$ws.Cells.Item(3, 5).Value = 89
$ws.Cells.Item(3, 6).Value = 2.01
Set-BlankCell 3 7

# Row 4 -- Portimonense SAD - SC União Torreense
$ws.Cells.Item(4, 1).Value = "Portimonense SAD - SC União Torreense : -:-'"
$ws.Cells.Item(4, 2).Value = "SC União Torreense"
$ws.Cells.Item(4, 3).Value = 29
Set-BlankCell 4 4
$ws.Cells.Item(4, 5).Value = 93
$ws.Cells.Item(4, 6).Value = 1.66
Set-BlankCell 4 7

# Row 5 -- CA Lanús - Club Universidad de Chile
$ws.Cells.Item(5, 1).Value = "CA Lanús  - Club Universidad de Chile: 21:00"
$ws.Cells.Item(5, 2).Value = "CA Lanús"
$ws.Cells.Item(5, 3).Value = 28
$ws.Cells.Item(5, 4).Value = 72
Set-BlankCell 5 5
$ws.Cells.Item(5, 6).Value = 2.88
Set-BlankCell 5 7

# Row 6 -- FC Lugano - FC Luzern
$ws.Cells.Item(6, 1).Value = "FC Lugano  - FC Luzern: -:-'"
$ws.Cells.Item(6, 2).Value = "FC Lugano"
$ws.Cells.Item(6, 3).Value = 25
$ws.Cells.Item(6, 4).Value = 75
Set-BlankCell 6 5
$ws.Cells.Item(6, 6).Value = 1.85
Set-BlankCell 6 7

# The blank-format source cell itself must stay untouched/empty.
$blankFormatSource.ClearContents()
